$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# 1) Mark the previous "第八周周四" block as fully completed: set all
#    "完成情况" (status) cells for that block to "完成".
# ------------------------------------------------------------------
$ws.Range("C94").Value = "完成"
$ws.Range("C95").Value = "完成"
$ws.Range("C96").Value = "完成"
$ws.Range("C97").Value = "完成"
$ws.Range("C98").Value = "完成"

# ------------------------------------------------------------------
# 2) Append a brand-new weekly block (rows 103-111) for
#    "日期：2018.10.25 第九周周一", re-using the same layout/formatting
#    as the previous block (rows 92-100).
# ------------------------------------------------------------------

# Pre-create the merged regions for the new block so that pasted
# formatting (border handling etc.) lines up the same way it does for
# the existing blocks.
$ws.Range("A103:D103").Merge()
$ws.Range("A110:D111").Merge()

# Copy the whole previous block's formatting onto the new block.
$ws.Range("A92:D100").Copy()
$ws.Range("A103").PasteSpecial(-4122)

# Title row.
$ws.Range("A103").Value = "日期：2018.10.25 第九周周一"

# Column header row.
$ws.Range("A104").Value = "组员"
$ws.Range("B104").Value = "计划内容"
$ws.Range("C104").Value = "完成情况"
$ws.Range("D104").Value = "备注"

# Data rows.
$ws.Range("A105").Value = "邱志鹏"
$ws.Range("B105").Value = "将组员的代码合并形成最新版本，完成个人信息设置的剩余界面"
$ws.Range("C105").Value = "未完成"

$ws.Range("A106").Value = "黄立根"
$ws.Range("B106").Value = "继续编写环信EaseUI的聊天界面,修改头像和昵称的显示"
$ws.Range("C106").Value = "未完成"

$ws.Range("A107").Value = "黄俊贤"
$ws.Range("B107").Value = "百度地图的点聚合和地图功能完善"
$ws.Range("C107").Value = "未完成"

$ws.Range("A108").Value = "李达波"
$ws.Range("B108").Value = "完成群聊剩余的界面，继续完成后台数据接口开发。"
$ws.Range("C108").Value = "未完成"

$ws.Range("A109").Value = "冯德志"
$ws.Range("B109").Value = "继续完成地图的路线规划、定位"
$ws.Range("C109").Value = "未完成"

# Summary row.
$ws.Range("A110").Value = "总结："

# ------------------------------------------------------------------
# 3) Update the view selection to match the edited location.
# ------------------------------------------------------------------
[void]$ws.Range("C98").Select()
